# Applies the scheduled-runner price/profit updates to the Leve price
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 85.566666
$ws.Range("I28").Value = 86.89655
$ws.Range("K28").Value = 86.89655
$ws.Range("M28").Value = 398.10345
$ws.Range("H64").Value = 3808.1667
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3969.8
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3969.8
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4465.8
$ws.Range("H67").Value = 3808.1667
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3969.8
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3969.8
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5685.8
$ws.Range("H74").Value = 3733.1667
$ws.Range("I74").Value = 2974.75
$ws.Range("K74").Value = 2974.75
$ws.Range("M74").Value = -2038.75
$ws.Range("H77").Value = 3733.1667
$ws.Range("I77").Value = 2974.75
$ws.Range("K77").Value = 14873.75
$ws.Range("M77").Value = -10193.75
$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880
$ws.Range("H132").Value = 857.75
$ws.Range("I132").Value = 857.75
$ws.Range("K132").Value = 2573.25
$ws.Range("M132").Value = -43.25
$ws.Range("H138").Value = 3574.3467
$ws.Range("J138").Value = 3254.2415
$ws.Range("L138").Value = 9762.7245
$ws.Range("N138").Value = -20042.7245

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 192483.8
$ws.Range("I2").Value = 278480
$ws.Range("J2").Value = 1381.1111
$ws.Range("K2").Value = 278480
$ws.Range("L2").Value = 1381.1111
$ws.Range("M2").Value = -278367
$ws.Range("N2").Value = -1607.1111
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384
$ws.Range("H32").Value = 19311.049
$ws.Range("I32").Value = 15719.77
$ws.Range("K32").Value = 15719.77
$ws.Range("M32").Value = -15432.77
$ws.Range("H45").Value = 2142.0715
$ws.Range("I45").Value = 2495.5
$ws.Range("J45").Value = 1877
$ws.Range("K45").Value = 2495.5
$ws.Range("L45").Value = 1877
$ws.Range("M45").Value = -2118.5
$ws.Range("N45").Value = -2631
$ws.Range("H74").Value = 1198.2727
$ws.Range("I74").Value = 879.5925999999999
$ws.Range("J74").Value = 2632.3333
$ws.Range("K74").Value = 879.5925999999999
$ws.Range("L74").Value = 2632.3333
$ws.Range("M74").Value = -5.592599999999948
$ws.Range("N74").Value = -4380.3333
$ws.Range("H77").Value = 1198.2727
$ws.Range("I77").Value = 879.5925999999999
$ws.Range("J77").Value = 2632.3333
$ws.Range("K77").Value = 4397.963
$ws.Range("L77").Value = 13161.6665
$ws.Range("M77").Value = -29.96299999999974
$ws.Range("N77").Value = -21897.6665
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("K102").Value = 2500
$ws.Range("M102").Value = -878
$ws.Range("H116").Value = 192483.8
$ws.Range("I116").Value = 278480
$ws.Range("J116").Value = 1381.1111
$ws.Range("K116").Value = 278480
$ws.Range("L116").Value = 1381.1111
$ws.Range("M116").Value = -276186
$ws.Range("N116").Value = -5969.1111

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 192483.8
$ws.Range("I3").Value = 278480
$ws.Range("J3").Value = 1381.1111
$ws.Range("K3").Value = 278480
$ws.Range("L3").Value = 1381.1111
$ws.Range("M3").Value = -278366
$ws.Range("N3").Value = -1609.1111
$ws.Range("H86").Value = 501498.5
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 1000997
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 1000997
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -1003243
$ws.Range("H89").Value = 501498.5
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 1000997
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 5004985
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -5016217
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H105").Value = 2629.2693
$ws.Range("I105").Value = 2265.0417
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 2265.0417
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -518.0417000000002
$ws.Range("N105").Value = -10494
$ws.Range("H107").Value = 2665
$ws.Range("I107").Value = 2627.7778
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2627.7778
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -707.7777999999998
$ws.Range("N107").Value = -6840
$ws.Range("H134").Value = 4002.6604
$ws.Range("I134").Value = 4082.2046
$ws.Range("K134").Value = 12246.6138
$ws.Range("M134").Value = -9711.613799999999
$ws.Range("M102").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4188.0713
$ws.Range("I31").Value = 2141.8
$ws.Range("K31").Value = 2141.8
$ws.Range("M31").Value = -1846.8
$ws.Range("H34").Value = 4188.0713
$ws.Range("I34").Value = 2141.8
$ws.Range("K34").Value = 2141.8
$ws.Range("M34").Value = -1939.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 179
$ws.Range("J23").Value = 198.33333
$ws.Range("L23").Value = 594.99999
$ws.Range("N23").Value = -1064.99999
$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H113").Value = 28833.75
$ws.Range("I113").Value = 111027
$ws.Range("J113").Value = 1436
$ws.Range("K113").Value = 333081
$ws.Range("L113").Value = 4308
$ws.Range("M113").Value = -330911
$ws.Range("N113").Value = -8648
$ws.Range("H137").Value = 4414.4546
$ws.Range("I137").Value = 2158.8572
$ws.Range("J137").Value = 8361.75
$ws.Range("K137").Value = 6476.571599999999
$ws.Range("L137").Value = 25085.25
$ws.Range("M137").Value = -1376.571599999999
$ws.Range("N137").Value = -35285.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 7693027.5
$ws.Range("J55").Value = 875
$ws.Range("L55").Value = 875
$ws.Range("N55").Value = -1221
$ws.Range("H61").Value = 3194.4285
$ws.Range("I61").Value = 2893.5833
$ws.Range("K61").Value = 2893.5833
$ws.Range("M61").Value = -2691.5833
$ws.Range("H68").Value = 2268
$ws.Range("I68").Value = 1975.2
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1975.2
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1226.2
$ws.Range("N68").Value = -4498
$ws.Range("H69").Value = 64999.832
$ws.Range("I69").Value = 81666.336
$ws.Range("J69").Value = 48333.332
$ws.Range("K69").Value = 81666.336
$ws.Range("L69").Value = 48333.332
$ws.Range("M69").Value = -80855.336
$ws.Range("N69").Value = -49955.332
$ws.Range("H71").Value = 2268
$ws.Range("I71").Value = 1975.2
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 9876
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -6132
$ws.Range("N71").Value = -22488
$ws.Range("H72").Value = 64999.832
$ws.Range("I72").Value = 81666.336
$ws.Range("J72").Value = 48333.332
$ws.Range("K72").Value = 244999.008
$ws.Range("L72").Value = 144999.996
$ws.Range("M72").Value = -240943.008
$ws.Range("N72").Value = -153111.996
$ws.Range("H82").Value = 2045
$ws.Range("I82").Value = 1949.5
$ws.Range("K82").Value = 1949.5
$ws.Range("M82").Value = -1588.5
$ws.Range("H85").Value = 2045
$ws.Range("I85").Value = 1949.5
$ws.Range("K85").Value = 1949.5
$ws.Range("M85").Value = -701.5
$ws.Range("H113").Value = 3194.4285
$ws.Range("I113").Value = 2893.5833
$ws.Range("K113").Value = 2893.5833
$ws.Range("M113").Value = -723.5832999999998
$ws.Range("H132").Value = 4372.4546
$ws.Range("I132").Value = 3629.28
$ws.Range("K132").Value = 10887.84
$ws.Range("M132").Value = -8357.84

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 673.86957
$ws.Range("I107").Value = 613.5909
$ws.Range("K107").Value = 1840.7727
$ws.Range("M107").Value = 79.22730000000001
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 2612.8262
$ws.Range("I126").Value = 2480.8096
$ws.Range("K126").Value = 7442.4288
$ws.Range("M126").Value = -4972.4288
$ws.Range("M113").ClearContents()
